$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D cells are treated as text, matching the original inlineStr formatting
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '29.368.11'
$ws.Range("E2").Value = '  +0.50%  '

$ws.Range("D3").Value = '1.871.44'
$ws.Range("E3").Value = '  +0.77%  '

$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  +0.17%  '

$ws.Range("D5").Value = '0.7153'
$ws.Range("E5").Value = '  +0.71%  '

$ws.Range("D6").Value = '238.96'
$ws.Range("E6").Value = '  +0.53%  '

$ws.Range("E7").Value = '  +0.13%  '

$ws.Range("B8").Value = 'LidoStakedEther'
$ws.Range("C8").Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range("D8").Value = '1.873.82'
$ws.Range("E8").Value = '  +0.30%  '

$ws.Range("B9").Value = 'Dogecoin'
$ws.Range("C9").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D9").Value = '0.07799'
$ws.Range("E9").Value = '  -4.46%  '

$ws.Range("B10").Value = 'Cardano'
$ws.Range("C10").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D10").Value = '0.3065'
$ws.Range("E10").Value = '  +0.88%  '

$ws.Range("B11").Value = 'Solana'
$ws.Range("C11").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D11").Value = '25.15'
$ws.Range("E11").Value = '  +8.34%  '

$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").Value = '0.08250'
$ws.Range("E12").Value = '  +0.63%  '

$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.859.01'
$ws.Range("E13").Value = '  -0.29%  '

$ws.Range("D14").Value = '5.235'
$ws.Range("E14").Value = '  +1.21%  '

$ws.Range("B15").Value = 'Polygon'
$ws.Range("C15").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D15").Value = '0.7213'
$ws.Range("E15").Value = '  +1.95%  '

$ws.Range("B16").Value = 'Litecoin'
$ws.Range("C16").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D16").Value = '90.53'
$ws.Range("E16").Value = '  +1.20%  '

$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").Value = '29.566.34'
$ws.Range("E17").Value = '  +1.02%  '

$ws.Range("B18").Value = 'Uniswap'
$ws.Range("C18").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D18").Value = '5.844'
$ws.Range("E18").Value = '  +0.99%  '

$ws.Range("B19").Value = 'ShibaInu'
$ws.Range("C19").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D19").Value = '0.000007866'
$ws.Range("E19").Value = '  -0.61%  '

$ws.Range("B20").Value = 'BitcoinCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D20").Value = '241.57'
$ws.Range("E20").Value = '  +1.59%  '

$ws.Range("B21").Value = 'Avalanche'
$ws.Range("C21").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D21").Value = '13.26'
$ws.Range("E21").Value = '  -0.62%  '

$ws.Range("B22").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C22").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D22").Value = '2.166.77'
$ws.Range("E22").Value = '  +2.48%  '

$ws.Range("B23").Value = 'Dai'
$ws.Range("C23").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("E23").Value = '  +0.02%  '

$ws.Range("B24").Value = 'BinanceUSD'
$ws.Range("C24").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D24").Value = '1.002'
$ws.Range("E24").Value = '  +0.16%  '

$ws.Range("B25").Value = 'Chainlink'
$ws.Range("C25").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D25").Value = '7.780'
$ws.Range("E25").Value = '  +4.87%  '

$ws.Range("B26").Value = 'Stellar'
$ws.Range("C26").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D26").Value = '0.1548'
$ws.Range("E26").Value = '  +5.95%  '

$ws.Range("B27").Value = 'Monero'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D27").Value = '163.39'
$ws.Range("E27").Value = '  +0.61%  '

$ws.Range("B28").Value = 'Cosmos'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D28").Value = '8.991'
$ws.Range("E28").Value = '  +0.39%  '

$ws.Range("B29").Value = 'EthereumClassic'
$ws.Range("C29").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D29").Value = '18.31'
$ws.Range("E29").Value = '  +1.24%  '

$ws.Range("B30").Value = 'LidoDAOToken'
$ws.Range("C30").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D30").Value = '1.927'
$ws.Range("E30").Value = '  -1.42%  '

$ws.Range("B31").Value = 'Toncoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D31").Value = '1.358'
$ws.Range("E31").Value = '  -4.92%  '

$ws.Range("B32").Value = 'PancakeSwap'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D32").Value = '1.480'
$ws.Range("E32").Value = '  +0.02%  '

$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").Value = '4.326'
$ws.Range("E33").Value = '  -1.69%  '

$ws.Range("B34").Value = 'InternetComputer(DFINITY)'
$ws.Range("C34").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D34").Value = '4.088'
$ws.Range("E34").Value = '  +1.66%  '

$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").Value = '0.05245'
$ws.Range("E35").Value = '  +0.44%  '

$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D36").Value = '1.197'
$ws.Range("E36").Value = '  +2.52%  '

$ws.Range("B37").Value = 'ImmutableX'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D37").Value = '0.7164'
$ws.Range("E37").Value = '  +1.23%  '

$ws.Range("B38").Value = 'Frax'
$ws.Range("C38").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D38").Value = '1.004'
$ws.Range("E38").Value = '  +0.46%  '

$ws.Range("B39").Value = 'HuobiToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D39").Value = '2.681'
$ws.Range("E39").Value = '  +0.17%  '

$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").Value = '0.01865'
$ws.Range("E40").Value = '  +0.42%  '

$ws.Range("B41").Value = 'MXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D41").Value = '2.721'
$ws.Range("E41").Value = '  -0.27%  '

$ws.Range("B42").Value = 'Maker'
$ws.Range("C42").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D42").Value = '1.179.61'
$ws.Range("E42").Value = '  +3.52%  '

$ws.Range("B43").Value = 'TrustWalletToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D43").Value = '0.9084'
$ws.Range("E43").Value = '  -1.53%  '

$ws.Range("B44").Value = 'FraxShare'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D44").Value = '6.018'
$ws.Range("E44").Value = '  +2.38%  '

$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").Value = '71.93'
$ws.Range("E45").Value = '  +2.62%  '

$ws.Range("B46").Value = 'TheSandbox'
$ws.Range("C46").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D46").Value = '0.4302'
$ws.Range("E46").Value = '  +0.46%  '

$ws.Range("B47").Value = 'PaxDollar'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D47").Value = '1.001'
$ws.Range("E47").Value = '  +0.18%  '

$ws.Range("B48").Value = 'Quant'
$ws.Range("C48").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D48").Value = '102.49'
$ws.Range("E48").Value = '  +0.15%  '

$ws.Range("B49").Value = 'Mantle'
$ws.Range("C49").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D49").Value = '0.5372'
$ws.Range("E49").Value = '  -0.52%  '

$ws.Range("B50").Value = 'RenderToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D50").Value = '1.763'
$ws.Range("E50").Value = '  -0.76%  '

$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").Value = '9.122'
$ws.Range("E51").Value = '  -0.73%  '
